$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-08-22 Tuesday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-08-23 Wednesday", 2)

# Update the division problems in the table (only the 5 data rows have content;
# the other rows are intentionally blank answer lines).
$t = $d.Tables.Item(1)

$values = @{
    1  = @("87÷3=", "18÷7=", "20÷2=", "59÷5=", "63÷7=")
    5  = @("58÷3=", "63÷9=", "57÷6=", "43÷5=", "45÷2=")
    9  = @("60÷6=", "11÷9=", "90÷7=", "62÷5=", "56÷3=")
    13 = @("30÷5=", "88÷3=", "88÷4=", "13÷2=", "26÷8=")
    17 = @("52÷4=", "57÷7=", "74÷7=", "39÷7=", "86÷9=")
}

foreach ($rowIndex in $values.Keys) {
    $cols = $values[$rowIndex]
    for ($c = 1; $c -le $cols.Length; $c++) {
        $cell = $t.Cell($rowIndex, $c)
        $cell.Range.Text = $cols[$c - 1]
    }
}
